# api supports results by mission id and result id
# Edits the "function_parameters" worksheet:
#  - inserts a new "param_index" column (C) between function_id and kind
#  - renames the "Octopus_Params" kind to "Tests_Params"
#  - updates rows so each function_id can have multiple indexed parameters
#  - adds two new parameter rows (9 and 10)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("function_parameters")

$xlCenter = -4108

# 1) Insert a new column before the current "kind" column (column C),
#    data in columns C:E shifts right to D:F and keeps formatting.
$ws.Columns(3).Insert()

# 2) Make sure the whole used area (including the rows we are about to add)
#    uses the same centered style as the rest of the table. Column A is only
#    used for rows 1-7 (rows 8-10 never had an A cell, keep it that way).
$ws.Range("A1:F7").HorizontalAlignment = $xlCenter
$ws.Range("A1:F7").VerticalAlignment = $xlCenter
$ws.Range("B8:F10").HorizontalAlignment = $xlCenter
$ws.Range("B8:F10").VerticalAlignment = $xlCenter

# 3) Re-write the data so that the new unique strings are first introduced
#    in this order: Tests_Params, __, param_index (matches target workbook).
$ws.Cells.Item(5,4).Value = "Tests_Params"
$ws.Cells.Item(5,5).Value = "__"
$ws.Cells.Item(1,3).Value = "param_index"

# 4) Fill in the rest of row 1 (header)
$ws.Cells.Item(1,1).Value = "id"
$ws.Cells.Item(1,2).Value = "function_id"
$ws.Cells.Item(1,4).Value = "kind"
$ws.Cells.Item(1,5).Value = "value"
$ws.Cells.Item(1,6).Value = "type"

# 5) Row 2 (function_id 1, param_index 1)
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(2,4).Value = "text"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = "int"

# 6) Row 3 (function_id 1, param_index 2)
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(3,3).Value = 2
$ws.Cells.Item(3,4).Value = "text"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = "int"

# 7) Row 4 (function_id 2, param_index 1)
$ws.Cells.Item(4,2).Value = 2
$ws.Cells.Item(4,3).Value = 1
$ws.Cells.Item(4,4).Value = "text"
$ws.Cells.Item(4,5).Value = 4
$ws.Cells.Item(4,6).Value = "int"

# 8) Row 5 (function_id 3, param_index 1) - Tests_Params/__ already set above
$ws.Cells.Item(5,2).Value = 3
$ws.Cells.Item(5,3).Value = 1
$ws.Cells.Item(5,6).Value = "DataFrame"

# 9) Row 6 (function_id 4, param_index 1)
$ws.Cells.Item(6,2).Value = 4
$ws.Cells.Item(6,3).Value = 1
$ws.Cells.Item(6,4).Value = "Sys_Params"
$ws.Cells.Item(6,5).Value = "__"
$ws.Cells.Item(6,6).Value = "DataFrame"

# 10) Row 7 (function_id 5, param_index 1)
$ws.Cells.Item(7,2).Value = 5
$ws.Cells.Item(7,3).Value = 1
$ws.Cells.Item(7,4).Value = "text"
$ws.Cells.Item(7,5).Value = 66
$ws.Cells.Item(7,6).Value = "string"

# 11) Row 8 (function_id 2, param_index 2)
$ws.Cells.Item(8,2).Value = 2
$ws.Cells.Item(8,3).Value = 2
$ws.Cells.Item(8,4).Value = "Sys_Params"
$ws.Cells.Item(8,5).Value = "__"
$ws.Cells.Item(8,6).Value = "DataFrame"

# 12) Row 9 (new) (function_id 3, param_index 2)
$ws.Cells.Item(9,2).Value = 3
$ws.Cells.Item(9,3).Value = 2
$ws.Cells.Item(9,4).Value = "Tests_Params"
$ws.Cells.Item(9,5).Value = "__"
$ws.Cells.Item(9,6).Value = "DataFrame"

# 13) Row 10 (new) (function_id 3, param_index 3)
$ws.Cells.Item(10,2).Value = 3
$ws.Cells.Item(10,3).Value = 3
$ws.Cells.Item(10,4).Value = "Sys_Params"
$ws.Cells.Item(10,5).Value = "__"
$ws.Cells.Item(10,6).Value = "DataFrame"

# 14) Column width for the newly inserted "param_index" column (columns D:F
#     already kept their original widths - 16.75 / 15.875 / 22.25 - from the
#     column insert above, so they do not need to be touched). The column is
#     sized to fit its "param_index" header text (best-fit style column).
$ws.Columns(3).ColumnWidth = 10.7

# 15) Match the saved selection in the target file
$ws.Activate()
$ws.Range("C2").Select()
